$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add new columns I (I0) and J (IF) -----------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the look of the existing header cells (B1:H1): bold font, thin box
# border, centered horizontally and top-aligned vertically.
$hdr = $ws.Range("I1:J1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.VerticalAlignment = -4160     # xlTop
$hdr.Borders.LineStyle = 1         # xlContinuous (thin box border)

# --- Data rows 2-34 --------------------------------------------------------
# Column I is a constant "1" and column J mirrors the existing "IP" value
# (column H) for every row except row 33, which carries its own values.
for ($r = 2; $r -le 34; $r++) {
    if ($r -eq 33) {
        $ws.Cells.Item($r, 9).Value = 5
        $ws.Cells.Item($r, 10).Value = 9
    } else {
        $ws.Cells.Item($r, 9).Value = 1
        $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 8).Value2
    }
}
